$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append below the existing data (rows 2-4), as rows 5-7.
# Columns not listed here (J, K, L, M, N, O, X, Z, AB, AC, AF, AH..AS, AU, AV)
# have no value in the source data and are left untouched.

$rows = @(
    @{
        Row = 5
        A = 112395253; B = 78657; C = "Ovaliderad"; D = "LC"; E = 229497
        F = "Korallblylav"; G = "Parmeliella triptophylla"; H = "(Ach.) Müll.Arg."
        P = "Påterud, Vrm"; Q = 333022; R = 6626740; S = 10
        T = "Värmland"; U = "Eda"; V = "Värmland"; W = "Järnskog"
        Y = "2023-09-19"; AA = "2023-09-19"
        AD = $false; AE = $false; AG = $false
        AW = "Jan Rees"; AX = "Jan Rees"
    },
    @{
        Row = 6
        A = 112395251; B = 93216; C = "Ovaliderad"; D = "LC"; E = 2810
        F = "Västlig hakmossa"; G = "Rhytidiadelphus loreus"; H = "(Hedw.) Warnst."
        P = "Påterud, Vrm"; Q = 333055; R = 6626785; S = 10
        T = "Värmland"; U = "Eda"; V = "Värmland"; W = "Järnskog"
        Y = "2023-09-19"; AA = "2023-09-19"
        AD = $false; AE = $false; AG = $false
        AW = "Jan Rees"; AX = "Jan Rees"
    },
    @{
        Row = 7
        A = 112395254; B = 89503; C = "Ovaliderad"; D = "LC"; E = 5447
        F = "Vedticka"; G = "Fuscoporia viticola"; H = "(Schwein.) Murrill"
        P = "Påterud, Vrm"; Q = 333021; R = 6626691; S = 10
        T = "Värmland"; U = "Eda"; V = "Värmland"; W = "Järnskog"
        Y = "2023-09-19"; AA = "2023-09-19"
        AD = $false; AE = $false; AG = $false
        AW = "Jan Rees"; AX = "Jan Rees"
    }
)

foreach ($rd in $rows) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value = $rd.A
    $ws.Cells.Item($r, 2).Value = $rd.B
    $ws.Cells.Item($r, 3).Value = $rd.C
    $ws.Cells.Item($r, 4).Value = $rd.D
    $ws.Cells.Item($r, 5).Value = $rd.E
    $ws.Cells.Item($r, 6).Value = $rd.F
    $ws.Cells.Item($r, 7).Value = $rd.G
    $ws.Cells.Item($r, 8).Value = $rd.H

    $ws.Cells.Item($r, 16).Value = $rd.P
    $ws.Cells.Item($r, 17).Value = $rd.Q
    $ws.Cells.Item($r, 18).Value = $rd.R
    $ws.Cells.Item($r, 19).Value = $rd.S
    $ws.Cells.Item($r, 20).Value = $rd.T
    $ws.Cells.Item($r, 21).Value = $rd.U
    $ws.Cells.Item($r, 22).Value = $rd.V
    $ws.Cells.Item($r, 23).Value = $rd.W

    # Force plain text (not auto-converted to a date serial) via leading apostrophe
    $ws.Cells.Item($r, 25).Value = "'" + $rd.Y
    $ws.Cells.Item($r, 27).Value = "'" + $rd.AA

    $ws.Cells.Item($r, 30).Value = $rd.AD
    $ws.Cells.Item($r, 31).Value = $rd.AE
    $ws.Cells.Item($r, 33).Value = $rd.AG

    $ws.Cells.Item($r, 49).Value = $rd.AW
    $ws.Cells.Item($r, 50).Value = $rd.AX
}
